$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the working-hours end time for row 64 (E64), from 21:00 to 19:30
$ws.Range("E64").Value = 19.5 / 24

# Recalculate dependent formulas (F64, G64, F66, F67, F68 all derive from E64)
$excel.Calculate()

# Update the active selection to reflect where the user left off editing
$ws.Range("E65").Select()
